$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.27812
$ws.Range("N2").Value = 3.83436
$ws.Range("O2").Value = 0.2125522080155258
$ws.Range("P2").Value = 0.2125522080155258
$ws.Range("Q2").Value = 0.1646069446
$ws.Range("R2").Value = 1.4814625014
$ws.Range("S2").Value = 0.2125522080155258
$ws.Range("T2").Value = 0.2125522080155258

# Row 3
$ws.Range("O3").Value = 0.4702887506191235
$ws.Range("P3").Value = 0.4702887506191235
$ws.Range("S3").Value = 0.4702887506191235
$ws.Range("T3").Value = 0.4702887506191235

# Row 4
$ws.Range("M4").Value = 1.907142333333333
$ws.Range("N4").Value = 5.721427
$ws.Range("O4").Value = 0.3171590413653506
$ws.Range("P4").Value = 0.3171590413653506
$ws.Range("Q4").Value = 0.2456176825394445
$ws.Range("R4").Value = 2.210559142855
$ws.Range("S4").Value = 0.3171590413653506
$ws.Range("T4").Value = 0.3171590413653506
